$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.262.37"
$ws.Range("E2").Value = "  +0.91%  "

$ws.Range("D3").Value = "1.911.08"
$ws.Range("E3").Value = "  +1.27%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "'321.37"
$ws.Range("E5").Value = "  -2.94%  "

$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").Value = "'0.4728"
$ws.Range("E7").Value = "  +2.78%  "

$ws.Range("D8").Value = "'0.4069"
$ws.Range("E8").Value = "  +0.22%  "

$ws.Range("D9").Value = "'0.08035"
$ws.Range("E9").Value = "  +0.73%  "

$ws.Range("D10").Value = "'1.002"
$ws.Range("E10").Value = "  +1.23%  "

$ws.Range("D11").Value = "'22.48"
$ws.Range("E11").Value = "  +3.75%  "

$ws.Range("D12").Value = "1.913.21"
$ws.Range("E12").Value = "  +0.88%  "

$ws.Range("D13").Value = "'5.893"
$ws.Range("E13").Value = "  -0.19%  "

$ws.Range("D14").Value = "'7.127"
$ws.Range("E14").Value = "  +1.10%  "

$ws.Range("D15").Value = "'89.71"
$ws.Range("E15").Value = "  +1.49%  "

$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.27%  "

$ws.Range("D17").Value = "'0.06638"
$ws.Range("E17").Value = "  +1.43%  "

$ws.Range("D18").Value = "'0.00001029"
$ws.Range("E18").Value = "  +0.11%  "

$ws.Range("E19").Value = "  +1.69%  "

$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("D21").Value = "29.280.17"
$ws.Range("E21").Value = "  +0.81%  "

$ws.Range("D22").Value = "'5.513"
$ws.Range("E22").Value = "  +2.05%  "

$ws.Range("D23").Value = "'11.46"
$ws.Range("E23").Value = "  +2.12%  "

$ws.Range("D25").Value = "2.146.63"
$ws.Range("E25").Value = "  +1.05%  "

$ws.Range("D26").Value = "'155.32"
$ws.Range("E26").Value = "  -0.69%  "

$ws.Range("D27").Value = "'19.75"
$ws.Range("E27").Value = "  +0.90%  "

$ws.Range("D28").Value = "'6.039"
$ws.Range("E28").Value = "  +11.70%  "

$ws.Range("D29").Value = "'2.106"
$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("D30").Value = "'117.32"
$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("D31").Value = "'1.068"
$ws.Range("E31").Value = "  +9.53%  "

$ws.Range("D32").Value = "'0.09535"
$ws.Range("E32").Value = "  +2.29%  "

$ws.Range("D33").Value = "'1.421"
$ws.Range("E33").Value = "  +1.24%  "

$ws.Range("D34").Value = "'3.545"
$ws.Range("E34").Value = "  -1.55%  "

$ws.Range("D35").Value = "'5.385"
$ws.Range("E35").Value = "  +2.16%  "

$ws.Range("D36").Value = "'0.06073"
$ws.Range("E36").Value = "  +0.54%  "

$ws.Range("D37").Value = "'0.02244"
$ws.Range("E37").Value = "  +1.16%  "

$ws.Range("D38").Value = "'8.241"
$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("D39").Value = "'1.171"
$ws.Range("E39").Value = "  -0.90%  "

$ws.Range("D40").Value = "'0.5853"
$ws.Range("E40").Value = "  +1.75%  "

$ws.Range("D41").Value = "'2.519"
$ws.Range("E41").Value = "  +12.48%  "

$ws.Range("D42").Value = "'0.1835"
$ws.Range("E42").Value = "  +1.09%  "

$ws.Range("D43").Value = "'10.09"
$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("D44").Value = "'0.07895"
$ws.Range("E44").Value = "  +2.25%  "

$ws.Range("D45").Value = "'1.273"
$ws.Range("E45").Value = "  +1.07%  "

$ws.Range("D46").Value = "'0.5526"
$ws.Range("E46").Value = "  +1.58%  "

$ws.Range("D47").Value = "'12.12"
$ws.Range("E47").Value = "  +1.31%  "

$ws.Range("D48").Value = "'1.922"
$ws.Range("E48").Value = "  +1.82%  "

$ws.Range("D49").Value = "'113.02"
$ws.Range("E49").Value = "  +1.90%  "

$ws.Range("D50").Value = "'44.38"
$ws.Range("E50").Value = "  -1.75%  "

$ws.Range("D51").Value = "'0.2927"
$ws.Range("E51").Value = "  +7.34%  "
